$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record is inserted as row 89, pushing the existing rows
# 89..137 down to 90..138 (dimension grows from A1:R137 to A1:R138).
$ws.Rows(89).Insert()

# Populate the newly inserted row with the new observation's data.
$ws.Range("A89").Value = 11
$ws.Range("B89").Value = "Vega Monumental Concepción"
$ws.Range("C89").Value = "Bíobío"
$ws.Range("D89").Value = 45176
$ws.Range("E89").Value = 8
$ws.Range("F89").Value = 100112037
$ws.Range("G89").Value = "Cebollín"
$ws.Range("H89").Value = "Sin especificar"
$ws.Range("I89").Value = "Primera"
$ws.Range("J89").Value = 60
$ws.Range("K89").Value = 4000
$ws.Range("L89").Value = 4200
$ws.Range("M89").Value = 4100
$ws.Range("N89").Value = '$/paquete 36 unidades'
$ws.Range("O89").Value = "Región Metropolitana"
$ws.Range("P89").Value = 114
$ws.Range("Q89").Value = 36
$ws.Range("R89").Value = "Hortaliza"
